$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + week-of dates) ---
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -76.666666666666
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 31
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 55
$ws.Range("L15").Value = 14.814814814814
$ws.Range("M15").Value = 63.157894736842
$ws.Range("N15").Value = -51.5625
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -70.588235294117
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = -11.304347826087
$ws.Range("L16").Value = -23.880597014925
$ws.Range("M16").Value = -61.363636363636
$ws.Range("N16").Value = -88.409090909090
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -30
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 312
$ws.Range("J17").Value = 314
$ws.Range("K17").Value = -0.636942675159
$ws.Range("L17").Value = -4.294478527607
$ws.Range("M17").Value = 15.985130111524
$ws.Range("N17").Value = -44.778761061946
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 116.666666666667
$ws.Range("I18").Value = 100
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 8.695652173913
$ws.Range("L18").Value = -13.043478260869
$ws.Range("M18").Value = -65.635738831615
$ws.Range("N18").Value = -88.151658767772
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -2.941176470588
$ws.Range("I19").Value = 343
$ws.Range("J19").Value = 265
$ws.Range("K19").Value = 29.433962264150
$ws.Range("L19").Value = 6.191950464396
$ws.Range("M19").Value = -28.541666666666
$ws.Range("N19").Value = -89.513910119229
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -15.384615384615
$ws.Range("I20").Value = 184
$ws.Range("J20").Value = 139
$ws.Range("K20").Value = 32.374100719424
$ws.Range("L20").Value = 15.723270440251
$ws.Range("M20").Value = -14.814814814814
$ws.Range("N20").Value = -86.050037907505
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = -19.047619047619
$ws.Range("I21").Value = 1079
$ws.Range("J21").Value = 955
$ws.Range("K21").Value = 12.984293193717
$ws.Range("L21").Value = -1.009174311926
$ws.Range("M21").Value = -30.566280566280
$ws.Range("N21").Value = -84.52602896888
$ws.Range("M23").Value = -60
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -14.814814814814
$ws.Range("F24").Value = 109
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = 25.287356321839
$ws.Range("I24").Value = 1011
$ws.Range("J24").Value = 726
$ws.Range("K24").Value = 39.256198347107
$ws.Range("L24").Value = 24.660912453760
$ws.Range("M24").Value = 26.059850374064
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -60
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -23.529411764705
$ws.Range("I25").Value = 407
$ws.Range("J25").Value = 374
$ws.Range("K25").Value = 8.823529411764
$ws.Range("L25").Value = 0.992555831265
$ws.Range("M25").Value = -36.206896551724
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 44
$ws.Range("J26").Value = 37
$ws.Range("K26").Value = 18.918918918918
$ws.Range("L26").Value = 12.820512820512
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 38
$ws.Range("K27").Value = 15.789473684210
$ws.Range("L27").Value = 109.52380952381
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 30
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -21.052631578947
$ws.Range("M28").Value = -37.5
$ws.Range("N28").Value = -74.576271186440
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 24
$ws.Range("K29").Value = -31.428571428571
$ws.Range("L29").Value = -22.580645161290
$ws.Range("M29").Value = -38.461538461538
$ws.Range("N29").Value = -77.142857142857

# --- Cells changing from numeric to text ("0" / "***.*" placeholders) ---
# Set value first (apostrophe forces literal text even for "0"), then copy number format
# from a reference text cell (A14, style s=14) so the resulting style id matches.
$ws.Range("C15").Value = "'0"
$ws.Range("C16").Value = "'0"
$ws.Range("G23").Value = "'0"
$ws.Range("H23").Value = "***.*"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "***.*"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "***.*"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)

# --- Cells changing from text back to numeric ---
$ws.Range("C28").Value = 2
$ws.Range("C29").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
